$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to text format so numeric-looking strings
# (e.g. "0.610", "1.00", "62.935.28") are preserved exactly as text,
# matching the original inlineStr storage instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.935.28'
$ws.Range("E2").Value = '  -4.31%  '
$ws.Range("D3").Value = '3.290.73'
$ws.Range("E3").Value = '  -6.39%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '540.38'
$ws.Range("E5").Value = '  -3.04%  '
$ws.Range("D6").Value = '169.88'
$ws.Range("E6").Value = '  -5.95%  '
$ws.Range("D7").Value = '0.610'
$ws.Range("E7").Value = '  -4.66%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = '3.277.78'
$ws.Range("E9").Value = '  -6.55%  '
$ws.Range("D10").Value = '0.605'
$ws.Range("E10").Value = '  -4.71%  '
$ws.Range("E11").Value = '  -1.71%  '
$ws.Range("D12").Value = '52.16'
$ws.Range("E12").Value = '  -3.64%  '
$ws.Range("D13").Value = '0.0000262'
$ws.Range("E13").Value = '  -3.90%  '
$ws.Range("D14").Value = '8.79'
$ws.Range("E14").Value = '  -5.53%  '
$ws.Range("D15").Value = '3.812.89'
$ws.Range("E15").Value = '  -6.35%  '
$ws.Range("D16").Value = '17.91'
$ws.Range("E16").Value = '  -3.67%  '
$ws.Range("E17").Value = '  -4.17%  '
$ws.Range("D18").Value = '3.291.65'
$ws.Range("E18").Value = '  -5.98%  '
$ws.Range("D19").Value = '11.54'
$ws.Range("E19").Value = '  -5.23%  '
$ws.Range("D20").Value = '62.887.75'
$ws.Range("E20").Value = '  -4.34%  '
$ws.Range("D21").Value = '0.964'
$ws.Range("E21").Value = '  -3.58%  '
$ws.Range("D22").Value = '411.89'
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("D23").Value = '4.37'
$ws.Range("E23").Value = '  +5.58%  '
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("D25").Value = '13.30'
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("D26").Value = '82.43'
$ws.Range("E26").Value = '  -4.87%  '
$ws.Range("D27").Value = '10.50'
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("E28").Value = '  -6.02%  '
$ws.Range("D29").Value = '8.49'
$ws.Range("E29").Value = '  -6.88%  '
$ws.Range("D30").Value = '28.81'
$ws.Range("E30").Value = '  -5.42%  '
$ws.Range("D31").Value = '6.29'
$ws.Range("E31").Value = '  -3.84%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '572.86'
$ws.Range("E32").Value = '  -5.91%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '11.24'
$ws.Range("E33").Value = '  -4.50%  '
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("D35").Value = '57.51'
$ws.Range("E35").Value = '  -3.81%  '
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").Value = '0.144'
$ws.Range("E37").Value = '  -1.88%  '
$ws.Range("D38").Value = '34.79'
$ws.Range("E38").Value = '  -7.46%  '
$ws.Range("E39").Value = '  +2.72%  '
$ws.Range("D40").Value = '0.0₃0729'
$ws.Range("E40").Value = '  -8.25%  '
$ws.Range("D41").Value = '0.360'
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("D42").Value = '3.094.30'
$ws.Range("E42").Value = '  -9.45%  '
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '3.24'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").Value = '2.74'
$ws.Range("E45").Value = '  -4.67%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0397'
$ws.Range("E46").Value = '  -4.64%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '2.39'
$ws.Range("E47").Value = '  -6.12%  '
$ws.Range("E48").Value = '  -4.18%  '
$ws.Range("E49").Value = '  -4.34%  '
$ws.Range("D50").Value = '132.34'
$ws.Range("E50").Value = '  -4.22%  '
$ws.Range("E51").Value = '  -6.51%  '
